$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.389.11'
$ws.Range('E2').Value = '  -1.36%  '
$ws.Range('D3').Value = '2.919.44'
$ws.Range('E3').Value = '  -2.91%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '376.05'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +6.31%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '102.51'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.27%  '
$ws.Range('E7').Value = '  -2.76%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.24%  '
$ws.Range('E9').Value = '  -4.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.90'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.13%  '
$ws.Range('E11').Value = '  -0.50%  '
$ws.Range('E12').Value = '  -2.29%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.29'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.71%  '
$ws.Range('D14').Value = '3.380.19'
$ws.Range('E14').Value = '  -2.80%  '
$ws.Range('E15').Value = '  -4.06%  '
$ws.Range('D16').Value = '2.913.29'
$ws.Range('E16').Value = '  -3.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.925'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -8.69%  '
$ws.Range('D18').Value = '51.313.44'
$ws.Range('E18').Value = '  -1.55%  '
$ws.Range('E19').Value = '  -0.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.33'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.98%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.91'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.48%  '
$ws.Range('E22').Value = '  -2.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.28'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '261.56'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.09%  '
$ws.Range('E25').Value = '  +1.19%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.169'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.42%  '
$ws.Range('E27').Value = '  -5.13%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '25.64'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.28'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.67%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.86'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +7.94%  '
$ws.Range('E32').Value = '  -4.23%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '9.77'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.18%  '
$ws.Range('E34').Value = '  -3.47%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '51.21'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.25%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '33.92'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.73%  '
$ws.Range('E37').Value = '  +0.32%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0421'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.09%  '
$ws.Range('E39').Value = '  -10.85%  '
$ws.Range('E40').Value = '  -3.55%  '
$ws.Range('E41').Value = '  -10.76%  '
$ws.Range('E42').Value = '  -7.80%  '
$ws.Range('E43').Value = '  -2.52%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '122.52'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.55%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.58'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.36%  '
$ws.Range('E46').Value = '  -2.70%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.268'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +10.72%  '
$ws.Range('D48').Value = '2.023.47'
$ws.Range('E48').Value = '  -4.62%  '
$ws.Range('E49').Value = '  -2.28%  '
$ws.Range('E50').Value = '  -5.19%  '
$ws.Range('D51').Value = '3.205.12'
$ws.Range('E51').Value = '  -2.89%  '
